$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 491.07144
$ws.Range("J6").Value = 1500
$ws.Range("L6").Value = 4500
$ws.Range("N6").Value = -4724

$ws.Range("H8").Value = 317
$ws.Range("I8").Value = 80.40000000000001
$ws.Range("J8").Value = 1500
$ws.Range("K8").Value = 241.2
$ws.Range("L8").Value = 4500
$ws.Range("M8").Value = -102.2
$ws.Range("N8").Value = -4778

$ws.Range("H33").Value = 222.55556
$ws.Range("I33").Value = 160.36
$ws.Range("K33").Value = 160.36
$ws.Range("M33").Value = 68.63999999999999

$ws.Range("H112").Value = 1281.0714
$ws.Range("I112").Value = 665
$ws.Range("J112").Value = 1449.091
$ws.Range("K112").Value = 1995
$ws.Range("L112").Value = 4347.272999999999
$ws.Range("M112").Value = -887
$ws.Range("N112").Value = -6563.272999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 5005000
$ws.Range("I11").Value = 5005000
$ws.Range("K11").Value = 5005000
$ws.Range("M11").Value = -5004856

$ws.Range("H132").Value = 2604.7917
$ws.Range("I132").Value = 2464.0527
$ws.Range("K132").Value = 7392.158100000001
$ws.Range("M132").Value = -4862.158100000001

$ws.Range("H134").Value = 150616.5
$ws.Range("J134").Value = 150616.5
$ws.Range("L134").Value = 150616.5
$ws.Range("N134").Value = -160756.5

$ws.Range("H135").Value = 49738.168
$ws.Range("J135").Value = 49738.168
$ws.Range("L135").Value = 49738.168
$ws.Range("N135").Value = -59878.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 8620
$ws.Range("J15").Value = 8620
$ws.Range("L15").Value = 8620
$ws.Range("N15").Value = -9074

$ws.Range("H35").Value = 16987
$ws.Range("J35").Value = 24974
$ws.Range("L35").Value = 24974
$ws.Range("N35").Value = -25594

$ws.Range("H82").Value = 16380.9
$ws.Range("I82").Value = 9986.666999999999
$ws.Range("J82").Value = 19121.285
$ws.Range("K82").Value = 9986.666999999999
$ws.Range("L82").Value = 19121.285
$ws.Range("M82").Value = -9603.666999999999
$ws.Range("N82").Value = -19887.285

$ws.Range("H85").Value = 16380.9
$ws.Range("I85").Value = 9986.666999999999
$ws.Range("J85").Value = 19121.285
$ws.Range("K85").Value = 9986.666999999999
$ws.Range("L85").Value = 19121.285
$ws.Range("M85").Value = -8660.666999999999
$ws.Range("N85").Value = -21773.285

$ws.Range("H86").Value = 1796.9375
$ws.Range("I86").Value = 1603.8889
$ws.Range("J86").Value = 2839.4
$ws.Range("K86").Value = 1603.8889
$ws.Range("L86").Value = 2839.4
$ws.Range("M86").Value = -480.8888999999999
$ws.Range("N86").Value = -5085.4

$ws.Range("H89").Value = 1796.9375
$ws.Range("I89").Value = 1603.8889
$ws.Range("J89").Value = 2839.4
$ws.Range("K89").Value = 8019.4445
$ws.Range("L89").Value = 14197
$ws.Range("M89").Value = -2403.4445
$ws.Range("N89").Value = -25429

$ws.Range("H134").Value = 38845.867
$ws.Range("I134").Value = 42512.89
$ws.Range("K134").Value = 127538.67
$ws.Range("M134").Value = -125003.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1750
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 1500
$ws.Range("M13").Value = -1861
$ws.Range("N13").Value = -1778

$ws.Range("H16").Value = 1294.5454
$ws.Range("I16").Value = 1292.5
$ws.Range("J16").Value = 1300
$ws.Range("K16").Value = 1292.5
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = -1005.5
$ws.Range("N16").Value = -1874

$ws.Range("H17").Value = 999
$ws.Range("I17").Value = 999
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 999
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -825
$ws.Range("N17").ClearContents()

$ws.Range("H23").Value = 2952.8096
$ws.Range("I23").Value = 2009
$ws.Range("K23").Value = 2009
$ws.Range("M23").Value = -1769

$ws.Range("H27").Value = 2952.8096
$ws.Range("I27").Value = 2009
$ws.Range("K27").Value = 2009
$ws.Range("M27").Value = -1817

$ws.Range("H50").Value = 6673.6
$ws.Range("J50").Value = 7092
$ws.Range("L50").Value = 7092
$ws.Range("N50").Value = -8342

$ws.Range("H110").Value = 69702
$ws.Range("J110").Value = 69702
$ws.Range("L110").Value = 69702
$ws.Range("N110").Value = -77882

$ws.Range("H113").Value = 1294.5454
$ws.Range("I113").Value = 1292.5
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 1292.5
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = 877.5
$ws.Range("N113").Value = -5640

$ws.Range("H122").Value = 2046.9333
$ws.Range("I122").Value = 1815.5454
$ws.Range("J122").Value = 2683.25
$ws.Range("K122").Value = 5446.6362
$ws.Range("L122").Value = 8049.75
$ws.Range("M122").Value = -2996.6362
$ws.Range("N122").Value = -12949.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 932.375
$ws.Range("J4").Value = 1520
$ws.Range("L4").Value = 4560
$ws.Range("N4").Value = -4784

$ws.Range("H5").Value = 809.5625
$ws.Range("I5").Value = 786.1818
$ws.Range("J5").Value = 861
$ws.Range("K5").Value = 2358.5454
$ws.Range("L5").Value = 2583
$ws.Range("M5").Value = -2246.5454
$ws.Range("N5").Value = -2807

$ws.Range("H31").Value = 1800
$ws.Range("J31").Value = 3100
$ws.Range("L31").Value = 9300
$ws.Range("N31").Value = -9876

$ws.Range("H131").Value = 1267821.1
$ws.Range("I131").Value = 3048.3333
$ws.Range("J131").Value = 1641032.8
$ws.Range("K131").Value = 9144.999899999999
$ws.Range("L131").Value = 4923098.4
$ws.Range("M131").Value = -4104.999899999999
$ws.Range("N131").Value = -4933178.4

$ws.Range("H135").Value = 809.5625
$ws.Range("I135").Value = 786.1818
$ws.Range("J135").Value = 861
$ws.Range("K135").Value = 7075.6362
$ws.Range("L135").Value = 7749
$ws.Range("M135").Value = -4540.6362
$ws.Range("N135").Value = -12819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H21").Value = 3000000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3000000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3000000
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -3000346

$ws.Range("H30").Value = 3000000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 3000000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 3000000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -3000210

$ws.Range("H43").Value = 17873.75
$ws.Range("J43").Value = 20355.715
$ws.Range("L43").Value = 20355.715
$ws.Range("N43").Value = -20657.715

$ws.Range("H70").Value = 5590.7
$ws.Range("I70").Value = 5585.2856
$ws.Range("K70").Value = 5585.2856
$ws.Range("M70").Value = -5315.2856

$ws.Range("H73").Value = 5590.7
$ws.Range("I73").Value = 5585.2856
$ws.Range("K73").Value = 5585.2856
$ws.Range("M73").Value = -4649.2856

$ws.Range("H80").Value = 172929.42
$ws.Range("I80").Value = 1600
$ws.Range("J80").Value = 201484.33
$ws.Range("K80").Value = 1600
$ws.Range("L80").Value = 201484.33
$ws.Range("M80").Value = -602
$ws.Range("N80").Value = -203480.33

$ws.Range("H83").Value = 172929.42
$ws.Range("I83").Value = 1600
$ws.Range("J83").Value = 201484.33
$ws.Range("K83").Value = 8000
$ws.Range("L83").Value = 1007421.65
$ws.Range("M83").Value = -3008
$ws.Range("N83").Value = -1017405.65

$ws.Range("H122").Value = 4390.4287
$ws.Range("I122").Value = 4622.1665
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 13866.4995
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -11416.4995
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 5003000
$ws.Range("I12").Value = 10000000
$ws.Range("J12").Value = 6000
$ws.Range("K12").Value = 10000000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = -9999830
$ws.Range("N12").Value = -6340

$ws.Range("H22").Value = 599.75
$ws.Range("J22").Value = 599.75
$ws.Range("L22").Value = 599.75
$ws.Range("N22").Value = -1189.75

$ws.Range("H27").Value = 599.75
$ws.Range("J27").Value = 599.75
$ws.Range("L27").Value = 599.75
$ws.Range("N27").Value = -813.75

$ws.Range("H46").Value = 1065.6666
$ws.Range("J46").Value = 647
$ws.Range("L46").Value = 647
$ws.Range("N46").Value = -1023

$ws.Range("H47").Value = 13099.1
$ws.Range("J47").Value = 13099.1
$ws.Range("L47").Value = 13099.1
$ws.Range("N47").Value = -14079.1

$ws.Range("H52").Value = 13099.1
$ws.Range("J52").Value = 13099.1
$ws.Range("L52").Value = 13099.1
$ws.Range("N52").Value = -13565.1

$ws.Range("H94").Value = 22330
$ws.Range("J94").Value = 22330
$ws.Range("L94").Value = 22330
$ws.Range("N94").Value = -23682

$ws.Range("H136").Value = 2115.5386
$ws.Range("I136").Value = 1087.75
$ws.Range("K136").Value = 3263.25
$ws.Range("M136").Value = -713.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 18910
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 18910
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 18910
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -19246

$ws.Range("H113").Value = 654
$ws.Range("I113").Value = 635.25
$ws.Range("J113").Value = 684
$ws.Range("K113").Value = 1905.75
$ws.Range("L113").Value = 2052
$ws.Range("M113").Value = 264.25
$ws.Range("N113").Value = -6392
